$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F, shifting the existing "District"/"NAME" column
# (old F) to G.
$ws.Columns.Item(6).Insert()

# New header for the inserted "Address" column.
$ws.Cells.Item(2, 6).Value = "Address"

# Row 3 sub-header cell stays blank (inlineStr with no text) for the new
# Address column.
$ws.Cells.Item(3, 6).Value = ""

# Fill in the address values (derived from the school name/address lines
# previously only embedded in column B) for each teacher row.
$addresses = @{
    4  = "Adarsha Vidyalaya Kollegal"
    5  = "Vijay High School KadangaArapattuMadikeri"
    6  = "Srimangala High SchoolSrimangala"
    7  = "Govt. High School PollibettaVirajpet"
    8  = "G H S GalibeeduMadikeri"
    9  = "G H S MaldareVirajpet"
    10 = "G H S BesurSomwarpet"
    11 = "Gonikoppal High School Gonikoppal"
    12 = "Rural High School PalyaKollegal"
    13 = "Janatha High SchoolHudikeri"
    14 = "Bharathi Vidya SamstheShanivarasanthe"
    15 = "G H P School ChikkathuruSomwarpet"
    16 = "S G M High School BhogainahundiGundalupet"
    17 = "G H SKesturYelandur"
    18 = "Parane High School ParaneMadikerei"
    19 = "G H S Chembu Madikeri"
    20 = "J P N High School Virajpet"
    21 = "Udaya High SchoolBettageriMadikeri"
    22 = "G HS Sagade"
    23 = "G M P SchoolSiddapurMadikeri"
    24 = "G H P S Nelliyahudi keriSomawarpet"
    25 = "JSS High SchoolGowdahalliYelandur"
    26 = "T G T G H P S ChenankoteVirjpet"
    27 = "G H S HeggalaVirajpet"
    28 = "J SS High School MudigundaKollegala"
    29 = "G H S MaralliHanur BlockKollegal"
    30 = "G H S MaddurYelandur"
    31 = "Govt High School NIREGundalpet"
    32 = "Govt. High School ShirangalaSomwarpet"
    33 = "G P U CollegeHigh School SectionMadikeri"
    34 = "G H S ChikkatuppurGundlupet"
    35 = "G H S MamballiYelandur"
    36 = "T S S S High School AgaraYelandur"
    37 = "G P U CollegeVenkataiahnachatra"
    38 = "K L High SchoolShanthalliSomwarpet"
    39 = "K b High School KuttandiVirajpet"
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item($row, 6).Value = $addresses[$row]
}
